# Data_exploration.pptx edit
#   1. Refresh the cached "datetimeFigureOut" field text (footer date) on the
#      slide master and every slide layout: 10/15/19 -> 10/11/21
#   2. Reword the first bullet on slide 2 ("Goals of data exploration")

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached date text wherever the auto "datetimeFigureOut" date
#    placeholder lives (the slide master + all custom/slide layouts).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "10/11/21"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout
}

# ---------------------------------------------------------------------------
# 2) Slide 2 ("Goals of data exploration") - first content bullet reworded.
# ---------------------------------------------------------------------------
$oldBullet = "Make sure data are correct" + " " + [char]0x2013 + " " + "should do alongside data wrangling"
$newBullet = "Make sure data are correct, especially after data wrangling"

$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*$oldBullet*") {
            [void]$tr.Replace($oldBullet, $newBullet)
        }
    }
}
